$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: "Images of People" feature gets a new highlighted note ---
# Add the note text first (matches the authoring order reflected by the
# shared-string table layout in the saved file).
$ws.Range("E32").Value = "Prepare but postpone application. Due to COVID, MSU is asking for iages of people with PPE"
$ws.Range("E32").Interior.Color = 49407
$ws.Range("E32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 29

# --- Row 18: "Explore Zoom timelines" task reworked into "Research format" ---
$ws.Range("B18").Value = "Research format"
$ws.Range("C18").Value = 1.25
$ws.Range("D18").Value = 0.75
$ws.Range("E18").Value = "A Zoom Timeline is too clunky"

# --- Update the view state: scroll down and move the active selection ---
$ws.Range("C19").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
